# Apply the Mon Mar 20 13:55:20 UTC 2023 GitHub Actions "cryptos list" price refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "1.150", "0.07400", "135.00")
# whose exact digits (incl. trailing zeros) must be preserved as TEXT, not coerced
# to a number. Force the text number format before assigning, then reset the cell
# style back to Normal so no stray formatting is left behind.
function Set-PriceText($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = 'Normal'
}

Set-PriceText 'D2' '28.086.31'
$ws.Range('E2').Value = '  +2.93%  '
Set-PriceText 'D3' '1.774.02'
$ws.Range('E3').Value = '  -0.54%  '
Set-PriceText 'D4' '1.008'
$ws.Range('E4').Value = '  +0.42%  '
Set-PriceText 'D5' '338.84'
$ws.Range('E5').Value = '  -0.56%  '
Set-PriceText 'D6' '1.006'
$ws.Range('E6').Value = '  +0.45%  '
Set-PriceText 'D7' '0.3828'
$ws.Range('E7').Value = '  -3.33%  '
Set-PriceText 'D8' '0.3421'
$ws.Range('E8').Value = '  -1.07%  '
Set-PriceText 'D9' '46.85'
$ws.Range('E9').Value = '  -2.23%  '
Set-PriceText 'D10' '1.150'
$ws.Range('E10').Value = '  -3.76%  '
Set-PriceText 'D11' '0.07400'
$ws.Range('E11').Value = '  -0.89%  '
Set-PriceText 'D12' '23.59'
$ws.Range('E12').Value = '  +8.63%  '
Set-PriceText 'D13' '1.004'
$ws.Range('E13').Value = '  +0.34%  '
Set-PriceText 'D14' '6.428'
$ws.Range('E14').Value = '  -0.85%  '
Set-PriceText 'D15' '7.324'
$ws.Range('E15').Value = '  +3.02%  '
Set-PriceText 'D16' '1.780.63'
$ws.Range('E16').Value = '  -0.17%  '
Set-PriceText 'D17' '0.00001078'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('E18').Value = '  -0.24%  '
Set-PriceText 'D19' '82.19'
$ws.Range('E19').Value = '  -2.50%  '
Set-PriceText 'D20' '1.005'
$ws.Range('E20').Value = '  +0.42%  '
Set-PriceText 'D21' '17.42'
$ws.Range('E21').Value = '  -1.28%  '
Set-PriceText 'D22' '6.409'
$ws.Range('E22').Value = '  -1.58%  '
Set-PriceText 'D23' '28.091.82'
$ws.Range('E23').Value = '  +2.94%  '
Set-PriceText 'D24' '12.08'
$ws.Range('E24').Value = '  -2.71%  '
Set-PriceText 'D25' '2.392'
$ws.Range('E25').Value = '  +0.50%  '
Set-PriceText 'D26' '20.73'
$ws.Range('E26').Value = '  -2.05%  '
Set-PriceText 'D27' '1.415'
$ws.Range('E27').Value = '  -3.99%  '
Set-PriceText 'D28' '2.409'
$ws.Range('E28').Value = '  -3.51%  '
Set-PriceText 'D29' '153.82'
$ws.Range('E29').Value = '  -2.34%  '
Set-PriceText 'D30' '1.982.90'
$ws.Range('E30').Value = '  -0.16%  '
Set-PriceText 'D31' '135.00'
$ws.Range('E31').Value = '  -1.09%  '
Set-PriceText 'D32' '4.028'
$ws.Range('E32').Value = '  +0.05%  '
Set-PriceText 'D33' '6.081'
$ws.Range('E33').Value = '  +1.99%  '
Set-PriceText 'D34' '0.08907'
$ws.Range('E34').Value = '  +0.93%  '
Set-PriceText 'D35' '12.75'
$ws.Range('E35').Value = '  -2.02%  '
Set-PriceText 'D36' '0.02402'
$ws.Range('E36').Value = '  -1.27%  '
Set-PriceText 'D37' '0.6815'
$ws.Range('E37').Value = '  -0.04%  '
Set-PriceText 'D38' '5.319'
$ws.Range('E38').Value = '  -1.70%  '
Set-PriceText 'D39' '0.06355'
$ws.Range('E39').Value = '  -1.63%  '
Set-PriceText 'D40' '0.2157'
$ws.Range('E40').Value = '  -2.14%  '
Set-PriceText 'D41' '1.244'
$ws.Range('E41').Value = '  -0.56%  '
Set-PriceText 'D42' '1.504'
$ws.Range('E42').Value = '  -7.19%  '
Set-PriceText 'D43' '8.285'
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-PriceText 'D44' '1.005'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-PriceText 'D45' '14.16'
$ws.Range('E45').Value = '  -1.62%  '
Set-PriceText 'D46' '0.6259'
$ws.Range('E46').Value = '  -2.08%  '
Set-PriceText 'D47' '3.860'
$ws.Range('E47').Value = '  -0.52%  '
Set-PriceText 'D48' '132.71'
$ws.Range('E48').Value = '  +0.25%  '
Set-PriceText 'D49' '2.065'
$ws.Range('E49').Value = '  -3.33%  '
Set-PriceText 'D50' '0.07515'
$ws.Range('E50').Value = '  +5.28%  '
$ws.Range('E51').Value = '  +3.28%  '
